$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect it so the cell values can be updated,
# then re-apply protection once the edits are complete.
$ws.Unprotect()

# Update the confidentiality / "as of" date note
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."
# Setting a multi-line value can trigger an automatic custom row height;
# restore the row to its default auto-fit height so it matches the source.
$ws.Rows.Item(7).AutoFit()

# Update the refreshed weight / percent-change figures for the holdings table
$ws.Range("D2").Value = 0.8465127595929526
$ws.Range("E2").Value = 0.008242454983515168

$ws.Range("D3").Value = 0.1534872404070474
$ws.Range("E3").Value = 0.009511376352107526

$ws.Range("E4").Value = 0.008437218222673737

# Restore sheet protection (content protected, formatting of rows/columns allowed)
$ws.Protect("D382", $true, $true, $true, $false, $true, $false, $false, $true, $true, $true, $true, $true, $true, $true, $true)
